# Auto-generated edit script applying cryptos.xlsx price/volume/name updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.0000", "0.3706")
# are preserved exactly as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.396.94'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '1.566.20'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '286.43'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').Value = '0.3706'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('D8').Value = '46.78'
$ws.Range('E8').Value = '  -3.69%  '
$ws.Range('D9').Value = '0.3287'
$ws.Range('E9').Value = '  -1.25%  '
$ws.Range('D10').Value = '1.148'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('D11').Value = '0.07435'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '20.50'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('D14').Value = '5.846'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').Value = '6.814'
$ws.Range('E15').Value = '  -1.28%  '
$ws.Range('D16').Value = '1.561.00'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '0.00001102'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '0.06696'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').Value = '86.40'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('D20').Value = '1.0000'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '6.325'
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').Value = '16.30'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').Value = '11.83'
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D24').Value = '22.411.83'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').Value = '2.326'
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('D26').Value = '2.575'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = '151.27'
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').Value = '19.40'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = '4.941'
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('D30').Value = '123.81'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '1.738.05'
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').Value = '1.057'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('D33').Value = '1.976'
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('D34').Value = '6.004'
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').Value = '9.691'
$ws.Range('E35').Value = '  -1.08%  '
$ws.Range('D36').Value = '0.08269'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '0.02404'
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '1.310'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.06334'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').Value = '0.2190'
$ws.Range('E40').Value = '  -2.08%  '
$ws.Range('D41').Value = '5.245'
$ws.Range('E41').Value = '  -2.18%  '
$ws.Range('D42').Value = '11.14'
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('D43').Value = '0.6139'
$ws.Range('E43').Value = '  -1.37%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '13.75'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.754'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5946'
$ws.Range('E47').Value = '  -0.73%  '
$ws.Range('D48').Value = '2.024'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').Value = '123.66'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '1.185'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('D51').Value = '0.07161'
$ws.Range('E51').Value = '  -0.54%  '
